$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5: switch the applied table style from the custom
#    "Table_0" style to the built-in "Themed Style 2 - Accent 1" style.
# ---------------------------------------------------------------------------
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTable) {
            $targetSlide = $sl
            $targetShape = $sh
        }
    }
}

if ($targetShape -ne $null) {
    $targetShape.Table.ApplyStyle("{0D97946A-5FF7-4D07-95F5-862F34F3BE28}")
}

# ---------------------------------------------------------------------------
# 2) Presentation theme: swap the "Integral / Red Violet" palette that is
#    currently applied to the deck for the standard "Office Theme" palette
#    (the two palettes that existed side-by-side in this file, one driving
#    the slides/master, the other only used by the notes master).
# ---------------------------------------------------------------------------
$colorScheme = $p.SlideMaster.ColorScheme
$colorScheme.Colors(1).RGB  = 0         # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477   # folHlink 954F72
